$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores values as text (inline strings) in the source file,
# even when the text looks numeric (e.g. "0.9992"). Pre-format just the cells
# receiving such numeric-looking text as Text ("@") so Excel keeps them as
# strings instead of silently coercing them into numeric cells.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D16", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.356.66"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "1.747.57"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "241.80"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "0.4818"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("D8").Value = "0.2617"
$ws.Range("E8").Value = "  +1.07%  "
$ws.Range("D9").Value = "0.06165"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "1.753.98"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "16.07"
$ws.Range("E11").Value = "  +3.12%  "
$ws.Range("D12").Value = "0.06941"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "0.6033"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "4.472"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "27.320.48"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "0.000007076"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("D20").Value = "11.48"
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Value = "1.960.11"
$ws.Range("E21").Value = "  +1.28%  "
$ws.Range("D22").Value = "4.442"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "8.436"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "5.115"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "142.15"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").Value = "15.23"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "1.839"
$ws.Range("E27").Value = "  +5.83%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "107.86"
$ws.Range("E28").Value = "  +2.56%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.384"
$ws.Range("E29").Value = "  -0.83%  "
$ws.Range("D30").Value = "3.954"
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").Value = "0.07975"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "3.671"
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "0.04673"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").Value = "2.600"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "1.013"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").Value = "0.6187"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").Value = "0.9244"
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").Value = "2.556"
$ws.Range("E38").Value = "  +6.77%  "
$ws.Range("D39").Value = "2.010"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").Value = "0.9995"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("D41").Value = "5.724"
$ws.Range("E41").Value = "  +5.62%  "
$ws.Range("D42").Value = "0.01494"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "99.80"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "0.3840"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "6.892"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "0.1155"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").Value = "7.879"
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("D49").Value = "29.86"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "1.248"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("D51").Value = "51.02"
$ws.Range("E51").Value = "  -0.37%  "

Write-Output "done"
